# Apply the "Changed scintillator temp sensor" edit:
# U3's Value and Device (TMP36GSZ) is replaced with TMP36FSZ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the temp sensor part number for U3 (row 23): Value (col B) and Device (col C)
$ws.Range("B23").Value = "TMP36FSZ"
$ws.Range("C23").Value = "TMP36FSZ"

# Update the selected cell shown in the sheet view
$ws.Range("C25").Select()
